# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates to the Carbuncle_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 33382.125
$ws.Range("I21").Value = 54764.25
$ws.Range("K21").Value = 54764.25
$ws.Range("M21").Value = -54296.25
$ws.Range("H23").Value = 33382.125
$ws.Range("I23").Value = 54764.25
$ws.Range("K23").Value = 54764.25
$ws.Range("M23").Value = -54530.25
$ws.Range("H87").Value = 17224.08
$ws.Range("J87").Value = 17224.08
$ws.Range("L87").Value = 17224.08
$ws.Range("N87").Value = -19720.08
$ws.Range("H90").Value = 17224.08
$ws.Range("J90").Value = 17224.08
$ws.Range("L90").Value = 51672.24000000001
$ws.Range("N90").Value = -64152.24000000001
$ws.Range("H112").Value = 1174.4546
$ws.Range("J112").Value = 1175.2333
$ws.Range("L112").Value = 3525.699900000001
$ws.Range("N112").Value = -5741.699900000001
$ws.Range("H113").Value = 11185.9375
$ws.Range("I113").Value = 3819
$ws.Range("J113").Value = 14534.546
$ws.Range("K113").Value = 3819
$ws.Range("L113").Value = 14534.546
$ws.Range("M113").Value = -565
$ws.Range("N113").Value = -21042.546
$ws.Range("H135").Value = 553.9524
$ws.Range("I135").Value = 543.3171
$ws.Range("K135").Value = 4889.8539
$ws.Range("M135").Value = -2354.8539
$ws.Range("H137").Value = 1004.75
$ws.Range("I137").Value = 953.5128
$ws.Range("J137").Value = 3003
$ws.Range("K137").Value = 2860.5384
$ws.Range("L137").Value = 9009
$ws.Range("M137").Value = -310.5383999999999
$ws.Range("N137").Value = -14109
$ws.Range("H138").Value = 3052.8765
$ws.Range("I138").Value = 1476.9302
$ws.Range("J138").Value = 4836.184
$ws.Range("K138").Value = 4430.7906
$ws.Range("L138").Value = 14508.552
$ws.Range("M138").Value = 709.2093999999997
$ws.Range("N138").Value = -24788.552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2880.73
$ws.Range("I32").Value = 2809.8284
$ws.Range("J32").Value = 9900
$ws.Range("K32").Value = 2809.8284
$ws.Range("L32").Value = 9900
$ws.Range("M32").Value = -2522.8284
$ws.Range("N32").Value = -10474
$ws.Range("H41").Value = 27381.875
$ws.Range("I41").Value = 9538.75
$ws.Range("J41").Value = 45225
$ws.Range("K41").Value = 9538.75
$ws.Range("L41").Value = 45225
$ws.Range("M41").Value = -9124.75
$ws.Range("N41").Value = -46053
$ws.Range("H61").Value = 740.05
$ws.Range("I61").Value = 740.05
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 740.05
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -528.05
$ws.Range("N61").ClearContents()
$ws.Range("H124").Value = 30885.334
$ws.Range("J124").Value = 30885.334
$ws.Range("L124").Value = 30885.334
$ws.Range("N124").Value = -40705.334
$ws.Range("H136").Value = 740.05
$ws.Range("I136").Value = 740.05
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2220.15
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 329.8500000000004
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1417.2391
$ws.Range("I134").Value = 813.75
$ws.Range("J134").Value = 3589.8
$ws.Range("K134").Value = 2441.25
$ws.Range("L134").Value = 10769.4
$ws.Range("M134").Value = 93.75
$ws.Range("N134").Value = -15839.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2490.3225
$ws.Range("I31").Value = 1703.8572
$ws.Range("J31").Value = 3138
$ws.Range("K31").Value = 1703.8572
$ws.Range("L31").Value = 3138
$ws.Range("M31").Value = -1408.8572
$ws.Range("N31").Value = -3728
$ws.Range("H34").Value = 2490.3225
$ws.Range("I34").Value = 1703.8572
$ws.Range("J34").Value = 3138
$ws.Range("K34").Value = 1703.8572
$ws.Range("L34").Value = 3138
$ws.Range("M34").Value = -1501.8572
$ws.Range("N34").Value = -3542
$ws.Range("H58").Value = 1130.8873
$ws.Range("I58").Value = 906.3390000000001
$ws.Range("J58").Value = 2234.9167
$ws.Range("K58").Value = 906.3390000000001
$ws.Range("L58").Value = 2234.9167
$ws.Range("M58").Value = -703.3390000000001
$ws.Range("N58").Value = -2640.9167
$ws.Range("H136").Value = 1130.8873
$ws.Range("I136").Value = 906.3390000000001
$ws.Range("J136").Value = 2234.9167
$ws.Range("K136").Value = 2719.017
$ws.Range("L136").Value = 6704.750100000001
$ws.Range("M136").Value = -169.0170000000003
$ws.Range("N136").Value = -11804.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3173.7292
$ws.Range("I131").Value = 464.44446
$ws.Range("J131").Value = 3798.9487
$ws.Range("K131").Value = 1393.33338
$ws.Range("L131").Value = 11396.8461
$ws.Range("M131").Value = 3646.66662
$ws.Range("N131").Value = -21476.8461
$ws.Range("H137").Value = 2080.074
$ws.Range("J137").Value = 2001.7368
$ws.Range("L137").Value = 6005.2104
$ws.Range("N137").Value = -16205.2104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3800
$ws.Range("I102").Value = 1200
$ws.Range("J102").Value = 4666.6665
$ws.Range("K102").Value = 1200
$ws.Range("L102").Value = 4666.6665
$ws.Range("M102").Value = 422
$ws.Range("N102").Value = -7910.6665
$ws.Range("H122").Value = 2505.394
$ws.Range("I122").Value = 2186.6365
$ws.Range("J122").Value = 3142.9092
$ws.Range("K122").Value = 6559.9095
$ws.Range("L122").Value = 9428.7276
$ws.Range("M122").Value = -4109.9095
$ws.Range("N122").Value = -14328.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2279.2856
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2279.2856
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6837.8568
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -11737.8568
$ws.Range("H132").Value = 7339.1577
$ws.Range("I132").Value = 9920.147999999999
$ws.Range("J132").Value = 5016.2666
$ws.Range("K132").Value = 29760.444
$ws.Range("L132").Value = 15048.7998
$ws.Range("M132").Value = -27230.444
$ws.Range("N132").Value = -20108.7998
$ws.Range("H138").Value = 46158.5
$ws.Range("J138").Value = 46158.5
$ws.Range("L138").Value = 46158.5
$ws.Range("N138").Value = -56438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 56751.75
$ws.Range("J18").Value = 56751.75
$ws.Range("L18").Value = 56751.75
$ws.Range("N18").Value = -57097.75
$ws.Range("H81").Value = 1456.0526
$ws.Range("I81").Value = 1081.2222
$ws.Range("J81").Value = 1793.4
$ws.Range("K81").Value = 2162.4444
$ws.Range("L81").Value = 3586.8
$ws.Range("M81").Value = -1101.4444
$ws.Range("N81").Value = -5708.8
$ws.Range("H84").Value = 1456.0526
$ws.Range("I84").Value = 1081.2222
$ws.Range("J84").Value = 1793.4
$ws.Range("K84").Value = 10812.222
$ws.Range("L84").Value = 17934
$ws.Range("M84").Value = -5508.222
$ws.Range("N84").Value = -28542
